$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refreshed prices & Volume(1h) percentages, and two row swaps
# (WrappedEther/Uniswap rows 17-18, and ApeXProtocol/Stellar rows 42-43).
# Cells that hold purely numeric-looking text (e.g. "122.40") are forced to stay
# text so Excel does not coerce them into numbers and drop trailing zeros.

$ws.Range('D2').Value = '63.195.50'
$ws.Range('E2').Value = '  -7.39%  '
$ws.Range('D3').Value = '3.540.85'
$ws.Range('E3').Value = '  -3.82%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '390.96'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '122.40'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.07%  '
$ws.Range('D7').Value = '3.528.82'
$ws.Range('E7').Value = '  -3.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.583'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -10.92%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  -12.26%  '
$ws.Range('E11').Value = '  -23.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000322'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -27.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.57'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.26%  '
$ws.Range('D14').Value = '4.100.64'
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.10'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.28%  '
$ws.Range('E16').Value = '  -2.94%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.549.34'
$ws.Range('E17').Value = '  -3.34%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.23'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.61'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.51%  '
$ws.Range('D20').Value = '63.242.32'
$ws.Range('E20').Value = '  -7.18%  '
$ws.Range('E21').Value = '  -9.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '391.03'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -15.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.85'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.30'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.77%  '
$ws.Range('E25').Value = '  -4.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.43'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.56'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.39%  '
$ws.Range('E28').Value = '  -8.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.63'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -14.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.65'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.74'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('E32').Value = '  -7.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.72'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.73%  '
$ws.Range('E34').Value = '  -5.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '36.32'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.44'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.86%  '
$ws.Range('E38').Value = '  -11.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = '0.0₃0644'
$ws.Range('E40').Value = '  -18.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.66'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.09'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +16.84%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.130'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -13.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '140.00'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.40'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +17.39%  '
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('E47').Value = '  -6.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.05'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.03%  '
$ws.Range('E49').Value = '  -9.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.62'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.98%  '
$ws.Range('E51').Value = '  -9.06%  '
